$wb = $excel.ActiveWorkbook

# 1. "Main Test": the test-run results in column A are refreshed to the
#    latest values produced by the script run (column D holds those values).
$wsMain = $wb.Worksheets.Item("Main Test")
for ($r = 2; $r -le 48; $r++) {
    $wsMain.Cells.Item($r, 1).Value2 = $wsMain.Cells.Item($r, 4).Value2
}

# 2. "Test_perl": fixed bug with std input - the perl example script now
#    flushes stdout and reads a value from stdin before processing.
$wsPerl = $wb.Worksheets.Item("Test_perl")
$perlScript = @'
use GD::Graph::points;
open(SRC, '<test_in.txt') or die $!;
my (@in1, @in2);
$|=1;
print "enter some value:\n";
my $val = <>;
print "value: $val\n";
while (<SRC>) {
    chomp;
    my @line = split /\t/;
    next if $line[1] eq "in2";
    print $line[1]." ... ".$line[2]."\n";
    push @in1, $line[1];
    push @in2, $line[2];
}
close SRC;
open(TGT, '>test_out.txt') or die $!;
print TGT $_."\n" for @in2;
close TGT;
print "plotting data\n";
my @data = ([@in1], [@in2]);
my $graph = GD::Graph::points->new(500, 300);
$graph->set(
                        x_label     => 'in1',
                        y_label     => 'in2',
                    ) or warn $mygraph->error;
my $gd = $graph->plot(\@data) or die $graph->error;
open(IMG, '>testdiagramperl.png') or die $!;
binmode IMG;
print IMG $gd->png;
print "finished \n";
'@
$wsPerl.Range("B3").Value2 = $perlScript

# 3. The workbook was left with the "Test_perl" sheet active/selected at B3.
$wsPerl.Activate()
$wsPerl.Range("B3").Select()

Write-Host "edit complete"
